# "error solve ifrs list"
#
# The IFRS financial-summary sheet ("company_list") had badly scaled /
# mis-pulled figures for the annual columns 2014..2018 (rows 2-6) and
# completely bogus rows for the *estimate* years 2019(E)-2021(E) (rows
# 7-9, all pulled from the wrong source block). This replaces the
# numeric data with the corrected figures and clears the cells that
# should no longer carry a value (including wiping the still-unreleased
# forecast rows back down to just their label columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : 2014/12 (IFRS연결) ----
$ws.Range("D2").Value  = 490
$ws.Range("E2").Value  = 87
$ws.Range("F2").Value  = 87
$ws.Range("G2").Value  = 46
$ws.Range("H2").Value  = -2
$ws.Range("I2").Value  = -2
$ws.Range("J2").Value  = 0
$ws.Range("K2").Value  = 855
$ws.Range("L2").Value  = 327
$ws.Range("M2").Value  = 528
$ws.Range("N2").Value  = 528
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value  = 46
$ws.Range("Q2").Value  = 52
$ws.Range("R2").Value  = 351
$ws.Range("S2").Value  = -403
$ws.Range("T2").Value  = 28
$ws.Range("U2").Value  = 24
$ws.Range("V2").Value  = 294
$ws.Range("W2").Value  = 17.85
$ws.Range("X2").Value  = -0.42
$ws.Range("Y2").Value  = -0.38
$ws.Range("Z2").Value  = -0.19
$ws.Range("AA2").Value = 62.03
$ws.Range("AB2").Value = 1048.39
$ws.Range("AC2").Value = -26
$ws.Range("AD2").Value = -286.33
$ws.Range("AE2").Value = 7154
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.69
$ws.Range("AI2").Value = -786.96
$ws.Range("AJ2").Value = 7378526

# ---- Row 3 : 2015/12 (IFRS연결) ----
$ws.Range("D3").Value  = 504
$ws.Range("E3").Value  = 111
$ws.Range("F3").Value  = 111
$ws.Range("G3").Value  = 100
$ws.Range("H3").Value  = 78
$ws.Range("I3").Value  = 78
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value  = 790
$ws.Range("L3").Value  = 211
$ws.Range("M3").Value  = 579
$ws.Range("N3").Value  = 579
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value  = 46
$ws.Range("Q3").Value  = 134
$ws.Range("R3").Value  = 29
$ws.Range("S3").Value  = -132
$ws.Range("T3").Value  = 19
$ws.Range("U3").Value  = 114
$ws.Range("V3").Value  = 186
$ws.Range("W3").Value  = 22.01
$ws.Range("X3").Value  = 15.44
$ws.Range("Y3").Value  = 14.06
$ws.Range("Z3").Value  = 9.46
$ws.Range("AA3").Value = 36.47
$ws.Range("AB3").Value = 1178.69
$ws.Range("AC3").Value = 1054
$ws.Range("AD3").Value = 8.31
$ws.Range("AE3").Value = 7953
$ws.Range("AF3").Value = 1.1
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 2.85
$ws.Range("AI3").Value = 23.39
$ws.Range("AJ3").Value = 7378526

# ---- Row 4 : 2016/12 (IFRS연결) ----
$ws.Range("D4").Value  = 489
$ws.Range("E4").Value  = 119
$ws.Range("F4").Value  = 119
$ws.Range("G4").Value  = 105
$ws.Range("H4").Value  = 81
$ws.Range("I4").Value  = 81
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value  = 723
$ws.Range("L4").Value  = 88
$ws.Range("M4").Value  = 636
$ws.Range("N4").Value  = 636
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value  = 46
$ws.Range("Q4").Value  = 122
$ws.Range("R4").Value  = 15
$ws.Range("S4").Value  = -155
$ws.Range("T4").Value  = 17
$ws.Range("U4").Value  = 105
$ws.Range("V4").Value  = 55
$ws.Range("W4").Value  = 24.29
$ws.Range("X4").Value  = 16.66
$ws.Range("Y4").Value  = 13.41
$ws.Range("Z4").Value  = 10.76
$ws.Range("AA4").Value = 13.77
$ws.Range("AB4").Value = 1314.52
$ws.Range("AC4").Value = 1104
$ws.Range("AD4").Value = 7.16
$ws.Range("AE4").Value = 8820
$ws.Range("AF4").Value = 0.9
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 3.8
$ws.Range("AI4").Value = 26.55
$ws.Range("AJ4").Value = 7378526

# ---- Row 5 : 2017/12 (IFRS연결) ----
$ws.Range("D5").Value  = 520
$ws.Range("E5").Value  = 81
$ws.Range("F5").Value  = 81
$ws.Range("G5").Value  = 66
$ws.Range("H5").Value  = 43
$ws.Range("I5").Value  = 67
$ws.Range("J5").Value  = -24
$ws.Range("K5").Value  = 5002
$ws.Range("L5").Value  = 3933
$ws.Range("M5").Value  = 1068
$ws.Range("N5").Value  = 659
$ws.Range("O5").Value  = 409
$ws.Range("P5").Value  = 46
$ws.Range("Q5").Value  = 33
$ws.Range("R5").Value  = -255
$ws.Range("S5").Value  = 356
$ws.Range("T5").Value  = 44
$ws.Range("U5").Value  = -10
$ws.Range("V5").Value  = 2535
$ws.Range("W5").Value  = 15.62
$ws.Range("X5").Value  = 8.29
$ws.Range("Y5").Value  = 10.42
$ws.Range("Z5").Value  = 1.5
$ws.Range("AA5").Value = 368.17
$ws.Range("AB5").Value = 1409.78
$ws.Range("AC5").Value = 914
$ws.Range("AD5").Value = 6.93
$ws.Range("AE5").Value = 9549
$ws.Range("AF5").Value = 0.66
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 5.52
$ws.Range("AI5").Value = 35.79
$ws.Range("AJ5").Value = 7378526

# ---- Row 6 : 2018/12 (IFRS연결) ----
$ws.Range("D6").Value  = 728
$ws.Range("E6").Value  = 12
$ws.Range("F6").Value  = 12
$ws.Range("G6").Value  = -99
$ws.Range("H6").Value  = -120
$ws.Range("I6").Value  = -6
$ws.Range("K6").Value  = 4402
$ws.Range("L6").Value  = 3486
$ws.Range("M6").Value  = 916
$ws.Range("N6").Value  = 623
$ws.Range("P6").Value  = 46
$ws.Range("Q6").Value  = -113
$ws.Range("R6").Value  = 403
$ws.Range("S6").Value  = -369
$ws.Range("T6").Value  = 63
$ws.Range("U6").Value  = -176
$ws.Range("V6").Value  = 2192
$ws.Range("W6").Value  = 1.69
$ws.Range("X6").Value  = -16.53
$ws.Range("Y6").Value  = -0.93
$ws.Range("Z6").Value  = -2.56
$ws.Range("AA6").Value = 380.62
$ws.Range("AB6").Value = 1471.69
$ws.Range("AC6").Value = -81
$ws.Range("AD6").Value = -72.94
$ws.Range("AE6").Value = 9285
$ws.Range("AF6").Value = 0.63
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 5.94
$ws.Range("AI6").Value = -394.07
$ws.Range("AJ6").Value = 7378526

# ---- Rows 7-9 : 2019(E)/2020(E)/2021(E) — forecast years have no
#      reliable source data yet, wipe the whole data block back to
#      just the row-number / label columns (A:C stay as-is). ----
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

Write-Output "ifrs list corrected (rows 2-9)"
